# Fruta / hortaliza, semanal
# Insert a new block of 3 rows (Especial / Primera / Segunda) for Mango at
# "Terminal La Palmera de La Serena" right before the existing row 350,
# pushing the former rows 350:461 down to 353:464.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 350; this shifts every row
# that was at 350:461 down to 353:464 (dimension grows from T461 to T464).
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(350).Insert()

# Row 350: Especial
$ws.Range("A350").Value = 8
$ws.Range("B350").Value = "Terminal La Palmera de La Serena"
$ws.Range("C350").Value = "Coquimbo"
$ws.Range("D350").Value = 44508
$ws.Range("D350").Style = $ws.Range("D353").Style
$ws.Range("E350").Value = 4
$ws.Range("F350").Value = "Fruta"
$ws.Range("G350").Value = 100108
$ws.Range("H350").Value = "Tropicales y subtropicales"
$ws.Range("I350").Value = 100108002
$ws.Range("J350").Value = "Mango"
$ws.Range("K350").Value = "Sin especificar"
$ws.Range("L350").Value = "Especial"
$ws.Range("M350").Value = 512
$ws.Range("N350").Value = 6500
$ws.Range("O350").Value = 7000
$ws.Range("P350").Value = 6750
$ws.Range("Q350").Value = "$/bandeja 4 kilos"
$ws.Range("R350").Value = "Perú"
$ws.Range("S350").Value = 1688
$ws.Range("T350").Value = 4

# Row 351: Primera
$ws.Range("A351").Value = 8
$ws.Range("B351").Value = "Terminal La Palmera de La Serena"
$ws.Range("C351").Value = "Coquimbo"
$ws.Range("D351").Value = 44508
$ws.Range("D351").Style = $ws.Range("D353").Style
$ws.Range("E351").Value = 4
$ws.Range("F351").Value = "Fruta"
$ws.Range("G351").Value = 100108
$ws.Range("H351").Value = "Tropicales y subtropicales"
$ws.Range("I351").Value = 100108002
$ws.Range("J351").Value = "Mango"
$ws.Range("K351").Value = "Sin especificar"
$ws.Range("L351").Value = "Primera"
$ws.Range("M351").Value = 512
$ws.Range("N351").Value = 6500
$ws.Range("O351").Value = 7000
$ws.Range("P351").Value = 6750
$ws.Range("Q351").Value = "$/bandeja 4 kilos"
$ws.Range("R351").Value = "Perú"
$ws.Range("S351").Value = 1688
$ws.Range("T351").Value = 4

# Row 352: Segunda
$ws.Range("A352").Value = 8
$ws.Range("B352").Value = "Terminal La Palmera de La Serena"
$ws.Range("C352").Value = "Coquimbo"
$ws.Range("D352").Value = 44508
$ws.Range("D352").Style = $ws.Range("D353").Style
$ws.Range("E352").Value = 4
$ws.Range("F352").Value = "Fruta"
$ws.Range("G352").Value = 100108
$ws.Range("H352").Value = "Tropicales y subtropicales"
$ws.Range("I352").Value = 100108002
$ws.Range("J352").Value = "Mango"
$ws.Range("K352").Value = "Sin especificar"
$ws.Range("L352").Value = "Segunda"
$ws.Range("M352").Value = 515
$ws.Range("N352").Value = 6500
$ws.Range("O352").Value = 7000
$ws.Range("P352").Value = 6751
$ws.Range("Q352").Value = "$/bandeja 4 kilos"
$ws.Range("R352").Value = "Perú"
$ws.Range("S352").Value = 1688
$ws.Range("T352").Value = 4
